# Kate Eisen - "Finished drafting test script; compiled data files out of R"
#
# The workbook holds one long/tidy-ish table of GC-MS peak areas: column A is
# a label column that mixes two kinds of names (the 5 sample/run names used
# as column headers in B1:F1, and the ~38 compound names used as row labels
# in A2:A38) with nothing distinguishing which is which once the sheet is
# fed into R. This script namespaces every label so the two kinds can be
# told apart downstream: sample names get a "Samp." prefix, compound names
# get a "Comp." prefix, and a handful of the compound names get their
# internal spacing normalised to dashes/no-space so they survive a round
# trip through R's data-frame column-name mangling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sample names used as the data-series headers (row 1, columns B:F) ---
$ws.Range("B1").Value = "Samp.rm3-1"
$ws.Range("C1").Value = "Samp.rm3-2"
$ws.Range("D1").Value = "Samp.rm3-4-1"
$ws.Range("E1").Value = "Samp.rm3-7"
$ws.Range("F1").Value = "Samp.rm3-8"

# --- Compound names, column A (row 1 is just the "Compound" header) ---
$ws.Range("A1").Value = "Compound"
$ws.Range("A2").Value = "Comp.2methylbutyronitrile"
$ws.Range("A3").Value = "Comp.3methylbutyronitrile"
$ws.Range("A4").Value = "Comp.b-myrcene"
$ws.Range("A5").Value = "Comp.cis-b-ocimene"
$ws.Range("A6").Value = "Comp.trans-b-ocimene"
$ws.Range("A7").Value = "Comp.nitro-2-methyl-butane"
$ws.Range("A8").Value = "Comp.nitro-3-methyl-butane"
$ws.Range("A9").Value = "Comp.isobutyronitrile1"
$ws.Range("A10").Value = "Comp.isobutyronitrile2"
$ws.Range("A11").Value = "Comp.furanoid-lin-oxide1"
$ws.Range("A12").Value = "Comp.furanoid-lin-oxide2"
$ws.Range("A13").Value = "Comp.pyran-lin-oxide-ketone"
$ws.Range("A14").Value = "Comp.2-methylbutyraldoxime1"
$ws.Range("A15").Value = "Comp.3-methylbutyraldoxime1"
$ws.Range("A16").Value = "Comp.2-methylbutyraldoxime2"
$ws.Range("A17").Value = "Comp.3-methylbutyraldoxime2"
$ws.Range("A18").Value = "Comp.linalool"
$ws.Range("A19").Value = "Comp.beta-caryophyllene"
$ws.Range("A20").Value = "Comp.beta-farnesene"
$ws.Range("A21").Value = "Comp.alpha-humulene"
$ws.Range("A22").Value = "Comp.alpha-terpineol"
$ws.Range("A23").Value = "Comp.Z-E-alpha-farnesene"
$ws.Range("A24").Value = "Comp.pyranoid-linalool-oxide1"
$ws.Range("A25").Value = "Comp.E-E-alpha-farnesene"
$ws.Range("A26").Value = "Comp.pyranoid-linalool-oxide2"
$ws.Range("A27").Value = "Comp.nerol"
$ws.Range("A28").Value = "Comp.geraniol"
$ws.Range("A29").Value = "Comp.2phenylethanol"
$ws.Range("A30").Value = "Comp.phenylacetonitrile"
$ws.Range("A31").Value = "Comp.farnesene epoxide1"
$ws.Range("A32").Value = "Comp.caryophyllene-oxide"
$ws.Range("A33").Value = "Comp.nerolidol"
$ws.Range("A34").Value = "Comp.farnesene-epoxide2"
$ws.Range("A35").Value = "Comp.nitrophenylethane"
$ws.Range("A36").Value = "Comp.phenylacetaldoxime"
$ws.Range("A37").Value = "Comp.isophytol"
$ws.Range("A38").Value = "Comp.farnesol"

# Retyping A31/A34/A35/A36 dropped their old highlight/fill formatting in the
# authored edit (the other cells on those rows keep theirs) - reset just the
# label cell's style back to the workbook default to match.
$ws.Range("A31").Style = "Normal"
$ws.Range("A34").Style = "Normal"
$ws.Range("A35").Style = "Normal"
$ws.Range("A36").Style = "Normal"

# Column A got noticeably wider (no longer auto "best fit") to fit the new
# "Comp."-prefixed names; ColumnWidth is in characters and Excel stores the
# column width with a fixed +5/6 character padding offset, so 27.667 here
# round-trips to the target stored width of 28.5.
$ws.Columns("A").ColumnWidth = 27.666666666666668

# Selection moved from the old active cell to the sample-name header row.
$null = $ws.Range("B1:F1").Select()
